# Sentia Log.xlsx - "Kleine updates aan logging en template."
#
# 1) Extend the "Traffic Manager" log entry (cell E13) with the extra
#    sentences about the Traffic Manager vs Load Balancer assumption.
# 2) Column D was switched from an auto best-fit width to a fixed,
#    narrower manual width (the long log texts are no longer fully
#    expanded in the grid).
# 3) The active selection ends up on E14 (the cell that was edited/
#    reviewed last).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = $ws.Range("E13").Value2
$addition = " Mijn aanname is om Traffic Manager te gebruiken en geen Load Balancer, omdat er in mijn omgeving geen sprake is van VM's die ge-loadbalanced moeten worden. Ik heb wel voor beide web-apps twee instances geconfigureerd zodat er interne load-balancing plaats vindt."
$ws.Range("E13").Value = $oldText + $addition

# Column D no longer auto-fits to the longest log entry; give it a
# fixed, compact manual width instead.
$ws.Columns.Item(4).ColumnWidth = 19.25

# Reflect the final cell the author left selected.
$ws.Range("E14").Select() | Out-Null
